$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used for the permuted data (1-indexed col numbers)
# D=4 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16 Q=17
$cols = @(4,8,9,10,11,12,13,14,15,16,17)

# Snapshot original values for rows 2-35 across the relevant columns
$snapshot = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowData
}

# Mapping: target row -> source row (values originally on source row move to target row)
$rowMap = @{
    2 = 17
    3 = 2
    4 = 3
    5 = 13
    6 = 26
    7 = 23
    8 = 32
    9 = 33
    10 = 34
    11 = 35
    12 = 21
    13 = 8
    14 = 12
    15 = 18
    16 = 31
    17 = 27
    18 = 28
    19 = 24
    20 = 19
    21 = 9
    22 = 10
    23 = 15
    24 = 16
    25 = 29
    26 = 30
    27 = 20
    28 = 7
    29 = 4
    30 = 14
    31 = 11
    32 = 6
    33 = 22
    34 = 25
    35 = 5
}

# Apply the permutation: write source row values into target row
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $srcData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcData[$c]
    }
}

Write-Host "Done applying row permutation."